$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new rule row at row 24 ("Set Response Due Date") ---
# Rows 24..33 (the "Set Due Date *" / "Set Owning Group" rule table) need to
# shift down by one to row 25..34, and the new rule occupies row 24.
#
# A plain Rows.Insert() pulls in stray new cell-style entries from this
# runtime's default "copy format from adjacent row" heuristic (row 23 above
# has very different borders/fills than row 24), so instead: first clone the
# formatting of the last table row into the new last row, then shift the
# cell values (not formats) up-to-down via a bottom-up copy, leaving row 24
# free for the new content.

# 1) Give row 34 the same formatting as row 33 (last existing table row) so it
#    is ready to receive the last row of shifted content.
$ws.Range("A33:F33").Copy()
$ws.Range("A34:F34").PasteSpecial(-4122)

# 2) Shift the values of rows 24-33 down into rows 25-34 (bottom-up so we
#    never overwrite a row before it's been read).
for ($r = 33; $r -ge 24; $r--) {
    for ($c = 1; $c -le 6; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r + 1, $c)
        $dst.Value2 = $src.Value2
    }
}

# 3) Row 24 keeps its original formatting (it already matched what the new
#    rule row needs), just replace its values with the new rule.
$ws.Range("B24").Value2 = "Set Response Due Date"
$ws.Range("C24").Value2 = "responseDueDate == null"
$ws.Range("D24").Value2 = "setResponseDueDate, java.time.LocalDate.now().plusDays(10)"
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()

# 4) The last row (now 34, was 33) used to carry the "setDueDate, null"
#    action in column D for the "Release Queue" rule; that action was
#    dropped, so clear it.
$ws.Range("D34").ClearContents()

$wb.Application.CutCopyMode = 0

# --- View-state refresh (scroll position / active selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
